$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: new checklist item - "submission is not valid, reject it and re-render "
$ws.Range("B11").Value = "submission is not valid, reject it and re-render "
$ws.Range("B11").Style = $ws.Range("B3").Style
$ws.Range("B11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 28.8

# Row 12: new checklist item - "feedback message should be next to the field "
$ws.Range("B12").Value = "feedback message should be next to the field "
$ws.Range("B12").Style = $ws.Range("B3").Style
$ws.Range("B12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 28.8

# Match the author's final selection (D12) recorded in the saved workbook view
$ws.Range("D12").Select()
